# early june 2025 update
# Updates the forecast date headers (row 1) and the VCI3M forecast values
# (rows 2-19, columns B:L) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> values for columns B..L (column B = index 2 ... column L = index 12)
$data = @{
    1  = @(45797, 45804, 45811, 45818, 45825, 45832, 45839, 45846, 45853, 45860, 45867)
    2  = @(67.59999999999999, 68.2, 68, 67, 65.09999999999999, 62.4, 59.1, 55.2, 50.9, 46.5, 42.3)
    3  = @(75.2, 75.90000000000001, 76.09999999999999, 75.7, 74.7, 72.90000000000001, 70.40000000000001, 67.2, 63.4, 59.1, 54.6)
    4  = @(65.90000000000001, 67.09999999999999, 67.40000000000001, 66.90000000000001, 65.5, 63.4, 60.7, 57.6, 54.1, 50.6, 47.1)
    5  = @(63.9, 64.3, 63.7, 62.2, 59.8, 56.8, 53.1, 49.1, 44.9, 40.7, 36.8)
    6  = @(68.09999999999999, 69.7, 70.90000000000001, 71.5, 71.5, 70.8, 69.40000000000001, 67.2, 64.5, 61.2, 57.7)
    7  = @(71.8, 70.5, 68.5, 66, 63, 59.6, 56, 52.4, 48.8, 45.4, 42.4)
    8  = @(54.6, 54.8, 54.7, 54.4, 53.7, 52.6, 51.1, 49.2, 46.9, 44.3, 41.5)
    9  = @(61.2, 59, 56.8, 54.7, 52.8, 51, 49.6, 48.4, 47.3, 46.5, 45.7)
    10 = @(56.8, 57.9, 58.1, 57.2, 55.4, 52.8, 49.3, 45.3, 41, 36.5, 32.1)
    11 = @(68.40000000000001, 69.8, 70.2, 69.59999999999999, 67.8, 64.90000000000001, 61, 56.3, 51, 45.6, 40.3)
    12 = @(80.59999999999999, 80.59999999999999, 79.40000000000001, 77, 73.5, 68.90000000000001, 63.6, 57.8, 51.8, 45.8, 40.3)
    13 = @(72.5, 74.2, 75.40000000000001, 76.09999999999999, 75.90000000000001, 74.7, 72.2, 68.5, 63.7, 57.9, 51.4)
    14 = @(64.09999999999999, 66.2, 67.8, 68.7, 68.5, 67.40000000000001, 65.2, 62.1, 58.3, 54.1, 49.9)
    15 = @(77.09999999999999, 73.40000000000001, 69, 64, 58.6, 53, 47.5, 42.3, 37.8, 34.1, 31.4)
    16 = @(56.9, 58.2, 58.6, 58.2, 56.9, 54.9, 52.2, 49, 45.4, 41.7, 38)
    17 = @(68.3, 70.2, 70.7, 69.7, 67.2, 63.2, 58.2, 52.2, 45.8, 39.3, 33.2)
    18 = @(74.40000000000001, 75.8, 75.7, 74.09999999999999, 71.09999999999999, 66.90000000000001, 61.7, 55.8, 49.5, 43.3, 37.5)
    19 = @(69.59999999999999, 72, 73.09999999999999, 72.7, 70.90000000000001, 67.7, 63.3, 57.8, 51.6, 45.1, 38.6)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # column B starts at index 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
